# Move " (Batista)" from the plain run " (Batista) + " into the preceding
# highlighted "SQL" run, so the highlighted run reads "SQL (Batista)" and
# the plain run that follows it is reduced to " + ".
#
#   "SQL"<highlight green> | " (Batista) + "<plain>
#   -> "SQL (Batista)"<highlight green> | " + "<plain>

$d = $word.ActiveDocument

$oldRun1Text = "SQL"
$oldRun2Text = " (Batista) + "
$newRun1Text = "SQL (Batista)"
$newRun2Text = " + "

# Locate the paragraph that still contains the original wording.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ("*" + $oldRun1Text + $oldRun2Text + "*")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Output "Target paragraph not found (already edited?)"
} else {
    $paraStart = $target.Range.Start
    $paraText = $target.Range.Text
    $offset = $paraText.IndexOf($oldRun1Text + $oldRun2Text)

    $run1Start = $paraStart + $offset
    $run1End = $run1Start + $oldRun1Text.Length

    $run1 = $d.Range($run1Start, $run1End)

    # Guard against re-applying the edit: if the character right after "SQL"
    # already shares its (green) highlight, the merge already happened.
    $boundaryCheck = $d.Range($run1End, $run1End + 1)
    $alreadyMerged = ($run1.Font.HighlightColorIndex -eq $boundaryCheck.Font.HighlightColorIndex) `
                     -and ($run1.Font.HighlightColorIndex -ne "")

    if ($alreadyMerged) {
        Write-Output "Already merged - skipping"
    } else {
        # Rewriting the text of the existing highlighted run extends it
        # in place, keeping its <w:rPr> (green highlight) intact.
        $run1.Text = $newRun1Text

        # The following run (formerly " (Batista) + ", plain formatting)
        # now starts right after the just-extended run and keeps the
        # same original length before being trimmed down to " + ".
        $run2Start = $run1.End
        $run2End = $run2Start + $oldRun2Text.Length
        $run2 = $d.Range($run2Start, $run2End)
        $run2.Text = $newRun2Text

        Write-Output "Updated run1=[$($run1.Text)] run2=[$($run2.Text)]"
    }
}
